$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format (style) from an existing formatted date cell (A2)
# down into the new date cells (A8:A17) so the numFmt/style is reused
# rather than creating a new custom style.
$ws.Range("A2").Copy()
$ws.Range("A8:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 8
$ws.Cells.Item(8, 1).Value = 44818
$ws.Cells.Item(8, 2).Value = "Tabaco"
$ws.Cells.Item(8, 3).Value = 5
$ws.Cells.Item(8, 4).Value = "Paquete"

# Row 9
$ws.Cells.Item(9, 1).Value = 44818
$ws.Cells.Item(9, 2).Value = "Copas"
$ws.Cells.Item(9, 3).Value = 10
$ws.Cells.Item(9, 4).Value = "Cervezas con Juanma"

# Row 10
$ws.Cells.Item(10, 1).Value = 44818
$ws.Cells.Item(10, 2).Value = "Tabaco"
$ws.Cells.Item(10, 3).Value = 28
$ws.Cells.Item(10, 4).Value = "Tabaco de liar"

# Row 11
$ws.Cells.Item(11, 1).Value = 44819
$ws.Cells.Item(11, 2).Value = "Cervezas"
$ws.Cells.Item(11, 3).Value = 7
$ws.Cells.Item(11, 4).Value = "Ecofamily"

# Row 12
$ws.Cells.Item(12, 1).Value = 44820
$ws.Cells.Item(12, 2).Value = "Cervezas"
$ws.Cells.Item(12, 3).Value = 7
$ws.Cells.Item(12, 4).Value = "Ecofamily"

# Row 13
$ws.Cells.Item(13, 1).Value = 44819
$ws.Cells.Item(13, 2).Value = "Cervezas"
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = "Ecofamily"

# Row 14
$ws.Cells.Item(14, 1).Value = 44821
$ws.Cells.Item(14, 2).Value = "Comida hostal"
$ws.Cells.Item(14, 3).Value = 12
$ws.Cells.Item(14, 4).Value = "Ecofamily"

# Row 15
$ws.Cells.Item(15, 1).Value = 44820
$ws.Cells.Item(15, 2).Value = "Chino"
$ws.Cells.Item(15, 3).Value = 4
$ws.Cells.Item(15, 4).Value = "Ecofamily"

# Row 16
$ws.Cells.Item(16, 1).Value = 44820
$ws.Cells.Item(16, 2).Value = "Kebak"
$ws.Cells.Item(16, 3).Value = 8
$ws.Cells.Item(16, 4).Value = "Ecofamily"

# Row 17
$ws.Cells.Item(17, 1).Value = 44820
$ws.Cells.Item(17, 2).Value = "Comida"
$ws.Cells.Item(17, 3).Value = 4
$ws.Cells.Item(17, 4).Value = "Ecofamily"

# Update selection to reflect the post-edit active cell (B18)
$ws.Range("B18").Select()
